# Imported arguments in workflows
# Adds a new "Invalid Scenario" test-data row to the TestData sheet and
# updates the active selection, mirroring the data the workflow now imports
# as arguments.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Duplicate the formatting (borders/style) of the last existing data row
# (row 3) into the new row 4, for columns A:D only (E stays empty).
$ws.Range("A3:D3").Copy($ws.Range("A4:D4"))

# Populate the new row with the "Invalid Scenario" test case data.
$ws.Range("A4").Value = "TC2-Return a sales slip - Invalid Scenario"
$ws.Range("B4").Value = 5649
$ws.Range("C4").Value = "Tilak"
$ws.Range("D4").Value = 740

# Match the author's final on-screen selection.
$ws.Range("F14").Select()
